# Updates the "cryptos" price/volume table (Sheet1) with refreshed values.
# D-column cells whose new text would otherwise be auto-parsed as a plain
# number (losing the trailing zero / thousands-dot formatting, e.g.
# "6.00" -> 6) are first forced to Text format ("@") so the literal
# string is preserved, matching the source data which stores these as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.066.76"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").Value = "2.251.74"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.31"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.27"
$ws.Range("E7").Value = "  +9.38%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +6.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.02"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.99"
$ws.Range("E12").Value = "  +4.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "2.587.51"
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.65"
$ws.Range("E15").Value = "  +4.65%  "
$ws.Range("D16").Value = "2.244.23"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.809"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "42.949.33"
$ws.Range("E18").Value = "  +4.45%  "
$ws.Range("E19").Value = "  +5.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.19"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  +5.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.28"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  +15.26%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.90"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.32"
$ws.Range("E28").Value = "  +28.57%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.64"
$ws.Range("E31").Value = "  +3.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.36"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0798"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("E34").Value = "  +4.88%  "
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("E36").Value = "  +8.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.36"
$ws.Range("E37").Value = "  +7.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0334"
$ws.Range("E38").Value = "  +17.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.99"
$ws.Range("E39").Value = "  +11.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.15"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.204"
$ws.Range("E42").Value = "  +7.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.32"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.88"
$ws.Range("E44").Value = "  +8.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.70"
$ws.Range("E45").Value = "  +5.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0996"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.463"
$ws.Range("E47").Value = "  +25.88%  "
$ws.Range("E48").Value = "  +8.87%  "
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("D51").Value = "2.461.46"
$ws.Range("E51").Value = "  +3.24%  "
